$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.286.87'
$ws.Range("E2").Value = '  -0.87%  '
$ws.Range("D3").Value = '1.704.95'
$ws.Range("E3").Value = '  -1.13%  '
$ws.Range("D4").Value = '1.003'
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("D5").Value = '223.54'
$ws.Range("E5").Value = '  -1.07%  '
$ws.Range("D6").Value = '0.5315'
$ws.Range("E6").Value = '  -0.94%  '
$ws.Range("E7").Value = '  -0.11%  '
$ws.Range("D8").Value = '0.2661'
$ws.Range("E8").Value = '  -0.54%  '
$ws.Range("D9").Value = '0.06596'
$ws.Range("E9").Value = '  -0.03%  '
$ws.Range("D10").Value = '20.78'
$ws.Range("E10").Value = '  -4.43%  '
$ws.Range("D11").Value = '0.07633'
$ws.Range("E11").Value = '  -1.37%  '
$ws.Range("D12").Value = '4.504'
$ws.Range("E12").Value = '  -2.59%  '
$ws.Range("D13").Value = '1.704.84'
$ws.Range("D14").Value = '1.936.79'
$ws.Range("E14").Value = '  -1.28%  '
$ws.Range("D15").Value = '0.5812'
$ws.Range("E15").Value = '  -0.94%  '
$ws.Range("D16").Value = '0.0₅8165'
$ws.Range("E16").Value = '  -1.70%  '
$ws.Range("D17").Value = '67.46'
$ws.Range("E17").Value = '  -0.85%  '
$ws.Range("D18").Value = '27.280.26'
$ws.Range("E18").Value = '  -0.96%  '
$ws.Range("D19").Value = '215.76'
$ws.Range("E19").Value = '  -3.23%  '
$ws.Range("E20").Value = '  -0.09%  '
$ws.Range("D21").Value = '4.629'
$ws.Range("E21").Value = '  -2.41%  '
$ws.Range("D22").Value = '10.35'
$ws.Range("E22").Value = '  -3.07%  '
$ws.Range("D23").Value = '5.972'
$ws.Range("E23").Value = '  -2.10%  '
$ws.Range("E24").Value = '  -0.16%  '
$ws.Range("D25").Value = '144.03'
$ws.Range("E25").Value = '  -2.65%  '
$ws.Range("D26").Value = '1.702'
$ws.Range("E26").Value = '  +0.53%  '
$ws.Range("D27").Value = '0.1201'
$ws.Range("E27").Value = '  -2.50%  '
$ws.Range("D28").Value = '7.211'
$ws.Range("E28").Value = '  -2.45%  '
$ws.Range("D29").Value = '16.16'
$ws.Range("E29").Value = '  -3.16%  '
$ws.Range("D30").Value = '0.05378'
$ws.Range("E30").Value = '  -2.86%  '
$ws.Range("D31").Value = '1.286'
$ws.Range("E31").Value = '  -1.31%  '
$ws.Range("D32").Value = '3.466'
$ws.Range("E32").Value = '  -2.27%  '
$ws.Range("D33").Value = '3.403'
$ws.Range("E33").Value = '  -1.76%  '
$ws.Range("D34").Value = '1.643'
$ws.Range("E34").Value = '  -1.15%  '
$ws.Range("D35").Value = '2.863'
$ws.Range("E35").Value = '  +1.46%  '
$ws.Range("D36").Value = '0.9476'
$ws.Range("E36").Value = '  -1.26%  '
$ws.Range("D37").Value = '2.409'
$ws.Range("E37").Value = '  -1.47%  '
$ws.Range("D38").Value = '0.5844'
$ws.Range("E38").Value = '  -1.64%  '
$ws.Range("D39").Value = '0.01636'
$ws.Range("E39").Value = '  -0.77%  '
$ws.Range("D40").Value = '5.805'
$ws.Range("E40").Value = '  -1.74%  '
$ws.Range("D41").Value = '1.003'
$ws.Range("E41").Value = '  -0.07%  '
$ws.Range("D42").Value = '1.041.54'
$ws.Range("E42").Value = '  -1.77%  '
$ws.Range("D43").Value = '0.8407'
$ws.Range("E43").Value = '  -1.85%  '
$ws.Range("D44").Value = '100.84'
$ws.Range("E44").Value = '  -0.83%  '
$ws.Range("D45").Value = '1.845.75'
$ws.Range("E45").Value = '  -1.20%  '
$ws.Range("D46").Value = '0.0₈112'
$ws.Range("E46").Value = '  -2.62%  '
$ws.Range("D47").Value = '57.77'
$ws.Range("E47").Value = '  -2.16%  '
$ws.Range("E48").Value = '  +1.73%  '
$ws.Range("D49").Value = '1.003'
$ws.Range("E49").Value = '  -0.24%  '
$ws.Range("D50").Value = '8.044'
$ws.Range("E50").Value = '  -2.04%  '
$ws.Range("D51").Value = '0.05225'
$ws.Range("E51").Value = '  -1.00%  '
